$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G holds strikeout count (K). Regenerated from source data.
$gValues = @{
    2 = 1
    3 = 1
    4 = 0
    5 = 2
    6 = 1
    7 = 3
    8 = 1
    9 = 2
    10 = 1
    11 = 0
    12 = 0
    13 = 2
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 3
    19 = 0
    20 = 1
    21 = 1
    22 = 3
    23 = 1
    24 = 3
    25 = 0
    26 = 1
    27 = 1
    28 = 0
    29 = 0
    30 = 1
    31 = 1
    32 = 2
    33 = 1
    34 = 1
    35 = 0
    36 = 1
    37 = 3
    38 = 2
    39 = 2
    41 = 0
    42 = 1
    43 = 0
    44 = 0
    45 = 0
    46 = 1
    47 = 1
    48 = 1
    49 = 2
    50 = 0
    51 = 1
    52 = 2
    53 = 1
    54 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value2 = $gValues[$row]
}

$wb.Save()
